$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3181
$ws1.Range("F5").Value = 2216
$ws1.Range("F8").Value = 1066
$ws1.Range("F9").Value = 1024
$ws1.Range("F16").Value = 7889
$ws1.Range("F18").Value = 2470
$ws1.Range("F22").Value = 461
$ws1.Range("F26").Value = 981
$ws1.Range("F27").Value = 1524
$ws1.Range("F28").Value = 11
$ws1.Range("F32").Value = 1908
$ws1.Range("F38").Value = 180

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3181
$ws4.Range("F7").Value = 2216
$ws4.Range("F10").Value = 1066
$ws4.Range("F12").Value = 1024
$ws4.Range("F19").Value = 7889
$ws4.Range("F21").Value = 2470
$ws4.Range("F26").Value = 461
$ws4.Range("F30").Value = 981
$ws4.Range("F31").Value = 1524
$ws4.Range("F32").Value = 11
$ws4.Range("F36").Value = 1908
$ws4.Range("F42").Value = 180
